$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 339
    $ws.Range("F4").Value = 499
    $ws.Range("G4").Value = 99
    $ws.Range("F5").Value = 4977
    $ws.Range("F9").Value = 753
}

# Row index for the "想去人数" column in row 10 differs between the two sheets:
# "展览" has the value in F10, "全部类型" has it (one row lower) in F11.
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F10").Value = 239

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F11").Value = 239
